$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete trailing rows (old rows 8-10, MuSCs-sending block)
$ws.Range("A8:A10").EntireRow.Delete() | Out-Null

# Refresh rows 2-7 with updated TPM-derived values
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Wnt11"
$ws.Range("C2").Value = "Fzd8"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 13.22273466666667
$ws.Range("H2").Value = 39.668204
$ws.Range("I2").Value = 0.9639959755780841
$ws.Range("J2").Value = 0.9639959755780841
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 3.083576666666666
$ws.Range("N2").Value = 9.250729999999999
$ws.Range("O2").Value = 0.2272509363535097
$ws.Range("P2").Value = 0.2272509363535097
$ws.Range("Q2").Value = 40.77331608765777
$ws.Range("R2").Value = 366.9598447889199
$ws.Range("S2").Value = 0.2190689880911346
$ws.Range("T2").Value = 0.2190689880911346

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Wnt11"
$ws.Range("C3").Value = "Fzd8"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 13.22273466666667
$ws.Range("H3").Value = 39.668204
$ws.Range("I3").Value = 0.9639959755780841
$ws.Range("J3").Value = 0.9639959755780841
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 6.453984666666667
$ws.Range("N3").Value = 19.361954
$ws.Range("O3").Value = 0.4756405360586227
$ws.Range("P3").Value = 0.4756405360586227
$ws.Range("Q3").Value = 85.33932679006844
$ws.Range("R3").Value = 768.0539411106159
$ws.Range("S3").Value = 0.4585155625823149
$ws.Range("T3").Value = 0.4585155625823149

$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Wnt11"
$ws.Range("C4").Value = "Fzd8"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 13.22273466666667
$ws.Range("H4").Value = 39.668204
$ws.Range("I4").Value = 0.9639959755780841
$ws.Range("J4").Value = 0.9639959755780841
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 4.031477000000001
$ws.Range("N4").Value = 12.094431
$ws.Range("O4").Value = 0.2971085275878677
$ws.Range("P4").Value = 0.2971085275878677
$ws.Range("Q4").Value = 53.30715068576934
$ws.Range("R4").Value = 479.764356171924
$ws.Range("S4").Value = 0.2864114249046346
$ws.Range("T4").Value = 0.2864114249046346

$ws.Range("A5").Value = "MuSCs"
$ws.Range("B5").Value = "Wnt11"
$ws.Range("C5").Value = "Fzd8"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.4938523333333333
$ws.Range("H5").Value = 1.481557
$ws.Range("I5").Value = 0.03600402442191584
$ws.Range("J5").Value = 0.03600402442191584
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.083576666666666
$ws.Range("N5").Value = 9.250729999999999
$ws.Range("O5").Value = 0.2272509363535097
$ws.Range("P5").Value = 0.2272509363535097
$ws.Range("Q5").Value = 1.522831531845555
$ws.Range("R5").Value = 13.70548378661
$ws.Range("S5").Value = 0.008181948262375003
$ws.Range("T5").Value = 0.008181948262375003

$ws.Range("A6").Value = "MuSCs"
$ws.Range("B6").Value = "Wnt11"
$ws.Range("C6").Value = "Fzd8"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.4938523333333333
$ws.Range("H6").Value = 1.481557
$ws.Range("I6").Value = 0.03600402442191584
$ws.Range("J6").Value = 0.03600402442191584
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.453984666666667
$ws.Range("N6").Value = 19.361954
$ws.Range("O6").Value = 0.4756405360586227
$ws.Range("P6").Value = 0.4756405360586227
$ws.Range("Q6").Value = 3.187315386930889
$ws.Range("R6").Value = 28.685838482378
$ws.Range("S6").Value = 0.01712497347630779
$ws.Range("T6").Value = 0.01712497347630779

$ws.Range("A7").Value = "MuSCs"
$ws.Range("B7").Value = "Wnt11"
$ws.Range("C7").Value = "Fzd8"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.4938523333333333
$ws.Range("H7").Value = 1.481557
$ws.Range("I7").Value = 0.03600402442191584
$ws.Range("J7").Value = 0.03600402442191584
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.031477000000001
$ws.Range("N7").Value = 12.094431
$ws.Range("O7").Value = 0.2971085275878677
$ws.Range("P7").Value = 0.2971085275878677
$ws.Range("Q7").Value = 1.990954323229667
$ws.Range("R7").Value = 17.918588909067
$ws.Range("S7").Value = 0.01069710268323304
$ws.Range("T7").Value = 0.01069710268323304

